$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-5
# from 45175 (2023-09-06) to 45183 (2023-09-14)
$ws.Range("C2:C5").Value = 45183
